$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The song previously listed as "Tsis Muaj Koj Pab (#148)" was renamed
$ws.Range("B4").Value = "Koj Tsis Pab, Kuv Ua Tsis Tau Dabtsi [Without Him] (#148)"

# Column B widens to fit the new, longer song title (best-fit style resize)
$ws.Columns.Item(2).ColumnWidth = 51.09

# Column C shrinks very slightly as part of the same autofit recalculation
$ws.Columns.Item(3).ColumnWidth = 140.8

# The active selection on the sheet moved to C4
[void]$ws.Range("C4").Select()
